$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.987.58"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "1.893.89"
$ws.Range("E3").Value = "  -2.00%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'0.7354"
$ws.Range("E5").Value = "  -2.09%  "

$ws.Range("D6").Value = "'242.76"
$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "'0.3084"
$ws.Range("E8").Value = "  -3.04%  "

$ws.Range("D9").Value = "'26.37"
$ws.Range("E9").Value = "  -4.23%  "

$ws.Range("D10").Value = "'0.06898"
$ws.Range("E10").Value = "  -1.20%  "

$ws.Range("D11").Value = "'0.7694"
$ws.Range("E11").Value = "  -1.34%  "

$ws.Range("D12").Value = "'0.07946"
$ws.Range("E12").Value = "  -0.58%  "

$ws.Range("D13").Value = "1.890.07"
$ws.Range("E13").Value = "  -2.22%  "

$ws.Range("D14").Value = "'5.219"
$ws.Range("E14").Value = "  -2.39%  "

$ws.Range("D15").Value = "'91.38"
$ws.Range("E15").Value = "  -3.09%  "

$ws.Range("D16").Value = "29.988.91"
$ws.Range("E16").Value = "  -1.15%  "

$ws.Range("D17").Value = "'14.09"
$ws.Range("E17").Value = "  -1.99%  "

$ws.Range("D18").Value = "'5.799"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007771"
$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'239.44"
$ws.Range("E20").Value = "  -5.30%  "

$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Value = "2.139.95"
$ws.Range("E22").Value = "  -2.27%  "

$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "'6.909"
$ws.Range("E24").Value = "  +3.55%  "

$ws.Range("D25").Value = "'9.295"
$ws.Range("E25").Value = "  -1.89%  "

$ws.Range("D26").Value = "'166.02"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("E28").Value = "  -4.37%  "

$ws.Range("D29").Value = "'2.024"
$ws.Range("E29").Value = "  -8.77%  "

$ws.Range("D30").Value = "'1.352"
$ws.Range("E30").Value = "  -0.69%  "

$ws.Range("D31").Value = "'1.536"
$ws.Range("E31").Value = "  +1.77%  "

$ws.Range("D32").Value = "'4.293"
$ws.Range("E32").Value = "  -1.77%  "

$ws.Range("D33").Value = "'4.052"
$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("D34").Value = "'0.05100"
$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("D35").Value = "'1.275"
$ws.Range("E35").Value = "  +0.38%  "

$ws.Range("E36").Value = "  -1.46%  "

$ws.Range("D37").Value = "'2.713"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").Value = "'0.01925"

$ws.Range("D39").Value = "'2.783"
$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("D40").Value = "'6.288"
$ws.Range("E40").Value = "  -1.65%  "

$ws.Range("D41").Value = "'74.22"
$ws.Range("E41").Value = "  -4.66%  "

$ws.Range("D42").Value = "'0.4442"
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").Value = "'1.934"
$ws.Range("E43").Value = "  -1.38%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "'0.8362"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("D46").Value = "'7.635"
$ws.Range("E46").Value = "  +2.40%  "

$ws.Range("D47").Value = "'101.03"
$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").Value = "'9.815"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").Value = "2.046.11"
$ws.Range("E49").Value = "  -1.92%  "

$ws.Range("D50").Value = "'36.41"
$ws.Range("E50").Value = "  -1.99%  "

$ws.Range("D51").Value = "'932.78"
$ws.Range("E51").Value = "  -5.21%  "
